$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60; this shifts existing rows 60-112 down to 61-113,
# preserving their data and formatting.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new data record.
$ws.Range("A60").Value2 = 5
$ws.Range("B60").Value2 = "Macroferia Regional de Talca"
$ws.Range("C60").Value2 = "Maule"
$ws.Range("D60").Value2 = 44879
$ws.Range("E60").Value2 = 7
$ws.Range("F60").Value2 = 100112022
$ws.Range("G60").Value2 = "Arveja Verde"
$ws.Range("H60").Value2 = "Sin especificar"
$ws.Range("I60").Value2 = "Primera"
$ws.Range("J60").Value2 = 500
$ws.Range("K60").Value2 = 16000
$ws.Range("L60").Value2 = 16000
$ws.Range("M60").Value2 = 16000
$ws.Range("N60").Value2 = "`$/saco 25 kilos"
$ws.Range("O60").Value2 = "Región del Maule"
$ws.Range("P60").Value2 = 640
$ws.Range("Q60").Value2 = 25
$ws.Range("R60").Value2 = "Hortaliza"
